$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 9: label + count of 1s per column (C..I), using rows 1-7
$ws.Range("A9").Value = "Cantidad de 0s ="
$ws.Range("C9").Formula = "=COUNTIF(C1:C7,1)"
$ws.Range("D9").Formula = "=COUNTIF(D1:D7,1)"
$ws.Range("E9").Formula = "=COUNTIF(E1:E7,1)"
$ws.Range("F9").Formula = "=COUNTIF(F1:F7,1)"
$ws.Range("G9").Formula = "=COUNTIF(G1:G7,1)"
$ws.Range("H9").Formula = "=COUNTIF(H1:H7,1)"
$ws.Range("I9").Formula = "=COUNTIF(I1:I7,1)"

# Row 10: label + count of 0s per column (C..I), using rows 1-7
$ws.Range("A10").Value = "Cantidad de 1s ="
$ws.Range("C10").Formula = "=COUNTIF(C1:C7,0)"
$ws.Range("D10").Formula = "=COUNTIF(D1:D7,0)"
$ws.Range("E10").Formula = "=COUNTIF(E1:E7,0)"
$ws.Range("F10").Formula = "=COUNTIF(F1:F7,0)"
$ws.Range("G10").Formula = "=COUNTIF(G1:G7,0)"
$ws.Range("H10").Formula = "=COUNTIF(H1:H7,0)"
$ws.Range("I10").Formula = "=COUNTIF(I1:I7,0)"
